# Update the date heading at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2026-02-11 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-12 Thursday", 2)

# Update the division-problem table. Cells are addressed by their
# (row, column) position so that duplicate problem text (e.g. "66÷7="
# appears twice, "82÷9=" appears twice) is replaced independently and
# correctly, rather than relying on document-wide text search/replace.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text  = "37÷7="
$t.Cell(1, 2).Range.Text  = "38÷6="
$t.Cell(1, 3).Range.Text  = "69÷2="
$t.Cell(1, 4).Range.Text  = "27÷3="
$t.Cell(1, 5).Range.Text  = "97÷3="

$t.Cell(5, 1).Range.Text  = "23÷5="
$t.Cell(5, 2).Range.Text  = "94÷9="
$t.Cell(5, 3).Range.Text  = "87÷5="
$t.Cell(5, 4).Range.Text  = "90÷4="
$t.Cell(5, 5).Range.Text  = "24÷3="

$t.Cell(9, 1).Range.Text  = "72÷7="
$t.Cell(9, 2).Range.Text  = "59÷8="
$t.Cell(9, 3).Range.Text  = "55÷4="
$t.Cell(9, 4).Range.Text  = "83÷7="
$t.Cell(9, 5).Range.Text  = "93÷7="

$t.Cell(13, 1).Range.Text = "74÷8="
$t.Cell(13, 2).Range.Text = "22÷4="
$t.Cell(13, 3).Range.Text = "12÷6="
$t.Cell(13, 4).Range.Text = "34÷5="
$t.Cell(13, 5).Range.Text = "90÷7="

$t.Cell(17, 1).Range.Text = "99÷7="
$t.Cell(17, 2).Range.Text = "10÷2="
$t.Cell(17, 3).Range.Text = "37÷5="
$t.Cell(17, 4).Range.Text = "60÷6="
$t.Cell(17, 5).Range.Text = "43÷5="
